$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(97).Insert()
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44494
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = 100112012
$ws.Range("G97").Value = "Espinaca"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 340
$ws.Range("K97").Value = 2300
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = 2394
$ws.Range("N97").Value = "$/docena de atados (3 kilos)"
$ws.Range("O97").Value = "Provincia de Quillota"
$ws.Range("P97").Value = 798
$ws.Range("Q97").Value = 3
$ws.Range("R97").Value = "Hortaliza"
